# repull data, push all data, mean calculation
# Update the dSF column (F) values for a set of rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 3
    5  = 0
    9  = -5
    10 = 2
    12 = -3
    16 = 2
    26 = -5
    36 = 0
    37 = -2
    41 = -2
    42 = 0
    49 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
